$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("IUCN 2019 Redlist")

# Insert a new column before column B (shifts Population estimate, Current
# Range, etc. one column to the right) to make room for species codes.
$ws.Columns.Item(2).EntireColumn.Insert()

# Header for the new column -- set first so the shared-string table gets
# the same ordering as a human typing the header before the data.
$ws.Range("B1").Value = "SpeciesCode"

# Species codes for each row (matches species listed in column A)
$ws.Range("B2").Value = "bw"
$ws.Range("B3").Value = "bp"
$ws.Range("B4").Value = "mn"
$ws.Range("B5").Value = "bs"
$ws.Range("B6").Value = "be"
$ws.Range("B7").Value = "ba"
$ws.Range("B8").Value = "bb"

# Re-apply the "Total removed" formula across the shifted column so it
# stays a single shared formula (E+F instead of the old D+E).
$ws.Range("G3:G8").Formula = "=E3+F3"

# Give the new column a sensible width (matches the other "short code"
# style columns on this sheet).
$ws.Columns.Item(2).ColumnWidth = 10.91

$ws.Range("C14").Select()
